$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.083.66'
$ws.Cells.Item(2, 5).Value = '  -0.15%  '
$ws.Cells.Item(3, 4).Value = '2.315.45'
$ws.Cells.Item(3, 5).Value = '  -0.05%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '302.31'
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  -0.39%  '
$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '98.97'
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  -1.97%  '
$ws.Cells.Item(7, 5).Value = '  +0.18%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.516'
$ws.Cells.Item(9, 4).Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  +1.25%  '
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '35.96'
$ws.Cells.Item(10, 4).Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  +3.38%  '
$origStyle = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0790'
$ws.Cells.Item(11, 4).Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  -1.02%  '
$ws.Cells.Item(12, 5).Value = '  -1.36%  '
$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '17.71'
$ws.Cells.Item(13, 4).Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  -2.48%  '
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.86'
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  -0.10%  '
$ws.Cells.Item(15, 4).Value = '2.677.02'
$ws.Cells.Item(15, 5).Value = '  -0.66%  '
$ws.Cells.Item(16, 4).Value = '2.301.68'
$ws.Cells.Item(16, 5).Value = '  -2.04%  '
$origStyle = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.792'
$ws.Cells.Item(17, 4).Style = $origStyle
$ws.Cells.Item(17, 5).Value = '  -3.67%  '
$ws.Cells.Item(18, 4).Value = '42.999.98'
$ws.Cells.Item(18, 5).Value = '  -0.27%  '
$ws.Cells.Item(19, 5).Value = '  +4.87%  '
$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.17'
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  -0.07%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0907'
$ws.Cells.Item(21, 5).Value = '  +0.03%  '
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '240.47'
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  +1.03%  '
$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.15'
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  -3.35%  '
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.45'
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  -1.29%  '
$ws.Cells.Item(26, 5).Value = '  -0.15%  '
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '25.09'
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  +0.74%  '
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '169.19'
$ws.Cells.Item(28, 4).Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  +0.53%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$origStyle = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.05'
$ws.Cells.Item(29, 4).Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  -2.20%  '
$ws.Cells.Item(30, 2).Value = 'Cosmos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '9.16'
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  -0.77%  '
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '33.35'
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  -2.89%  '
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.91'
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  +3.93%  '
$origStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.16'
$ws.Cells.Item(33, 4).Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  +2.30%  '
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  -0.05%  '
$origStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '18.29'
$ws.Cells.Item(35, 4).Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  +5.89%  '
$ws.Cells.Item(36, 5).Value = '  -0.76%  '
$origStyle = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.0692'
$ws.Cells.Item(37, 4).Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  -0.48%  '
$ws.Cells.Item(38, 5).Value = '  +0.64%  '
$ws.Cells.Item(39, 5).Value = '  -0.44%  '
$ws.Cells.Item(40, 5).Value = '  -2.82%  '
$origStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.109'
$ws.Cells.Item(41, 4).Style = $origStyle
$ws.Cells.Item(41, 5).Value = '  -0.57%  '
$ws.Cells.Item(42, 4).Value = '1.994.67'
$ws.Cells.Item(42, 5).Value = '  -0.66%  '
$ws.Cells.Item(43, 5).Value = '  +0.23%  '
$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '10.14'
$ws.Cells.Item(44, 4).Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  -0.33%  '
$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '17.37'
$ws.Cells.Item(45, 4).Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  -1.90%  '
$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.83'
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  -1.54%  '
$ws.Cells.Item(47, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.00'
$ws.Cells.Item(47, 4).Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  -14.94%  '
$ws.Cells.Item(48, 2).Value = 'BitcoinSV'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '76.39'
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  +8.35%  '
$ws.Cells.Item(49, 2).Value = 'MultiversX'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '54.57'
$ws.Cells.Item(49, 4).Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  -2.68%  '
$ws.Cells.Item(50, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(50, 4).Value = '2.544.19'
$ws.Cells.Item(50, 5).Value = '  +0.69%  '
$ws.Cells.Item(51, 2).Value = 'Stacks'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.54'
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  -0.46%  '
